# Edit: insert two new weekly price records (Primera/Segunda) for
# "Femacal de La Calera - Brócoli" dated 2022-07-27 (serial 44769),
# shifting the existing rows 690-751 down to 692-753.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 690 (old rows 690-751 move to 692-753)
$ws.Range("A690:A691").EntireRow.Insert()

# New row 690 - "Primera" quality
$ws.Cells.Item(690, 1).Value2 = 3
$ws.Cells.Item(690, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(690, 3).Value2 = "Coquimbo"
$ws.Cells.Item(690, 4).Value2 = 44769
$ws.Cells.Item(690, 5).Value2 = 5
$ws.Cells.Item(690, 6).Value2 = 100112023
$ws.Cells.Item(690, 7).Value2 = "Brócoli"
$ws.Cells.Item(690, 8).Value2 = "Sin especificar"
$ws.Cells.Item(690, 9).Value2 = "Primera"
$ws.Cells.Item(690, 10).Value2 = 1500
$ws.Cells.Item(690, 11).Value2 = 900
$ws.Cells.Item(690, 12).Value2 = 900
$ws.Cells.Item(690, 13).Value2 = 900
$ws.Cells.Item(690, 14).Value2 = "`$/unidad"
$ws.Cells.Item(690, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(690, 16).Value2 = 900
$ws.Cells.Item(690, 17).Value2 = 1
$ws.Cells.Item(690, 18).Value2 = "Hortaliza"

# New row 691 - "Segunda" quality
$ws.Cells.Item(691, 1).Value2 = 3
$ws.Cells.Item(691, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(691, 3).Value2 = "Coquimbo"
$ws.Cells.Item(691, 4).Value2 = 44769
$ws.Cells.Item(691, 5).Value2 = 5
$ws.Cells.Item(691, 6).Value2 = 100112023
$ws.Cells.Item(691, 7).Value2 = "Brócoli"
$ws.Cells.Item(691, 8).Value2 = "Sin especificar"
$ws.Cells.Item(691, 9).Value2 = "Segunda"
$ws.Cells.Item(691, 10).Value2 = 1400
$ws.Cells.Item(691, 11).Value2 = 700
$ws.Cells.Item(691, 12).Value2 = 700
$ws.Cells.Item(691, 13).Value2 = 700
$ws.Cells.Item(691, 14).Value2 = "`$/unidad"
$ws.Cells.Item(691, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(691, 16).Value2 = 700
$ws.Cells.Item(691, 17).Value2 = 1
$ws.Cells.Item(691, 18).Value2 = "Hortaliza"
